$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The queries table gains a new "Group" column (all rows tagged "A" for
# now, since MS2 results are split by collision method/energy downstream).
# Insert a new column in front of the existing "Formula" column (C) -
# everything from C onward (Formula, Monoisotopic, ion_mode, TOLERANCEPPM,
# RTMIN, RTMAX, QC_threshold) shifts one column to the right.
$ws.Range("C1").EntireColumn.Insert()

# The inserted column should carry the same formatting the "Formula"
# column had (it's now shifted to D), so copy that formatting over.
$ws.Range("D1:D7").Copy()
$ws.Range("C1:C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Give the new column a sensible width, matching its sibling text columns.
$ws.Range("C1").EntireColumn.ColumnWidth = 20.6197917

# Header + data
$ws.Range("C1").Value = "Group"
$ws.Range("C2:C7").Value = "A"

# Move the selection to mirror the post-edit state.
$ws.Range("D20").Select()
